$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original Datasets-tab query (currently in B2) moves down to B3 ("ProjectsTab" row),
# and a new, trimmed-down query (without Participants Count, with title REPLACE) takes its
# place in B2 ("DatasetsTab" row).

$oldQuery = $ws.Range("B2").Value2

$newQuery = "SELECT DISTINCT`n    REPLACE(ds.dataset_title, '  ', ' ') AS ""Title"", `n    ds.dataset_source_id AS ""Source ID"", `n    ds.primary_disease AS ""Primary Disease"",`n    -- CAST(ds.participant_count AS INT) AS ""Participants Count"",`n    CAST(ds.sample_count AS INT) AS ""Sample Count""`nFROM df_cedcd ds`nORDER BY ds.dataset_title ASC;"

$ws.Range("B3").Value = $oldQuery
$ws.Range("B2").Value = $newQuery

$ws.Range("C2").Select()
